# Update "Lương" sheet: rows/labels were re-numbered (3 new rows inserted per
# location: "Tổng công" / "Phụ cấp" / "Lương công tác"), totals recomputed for
# 26 worked days at SÓC TRĂNG (was 24), plus three new "Tổng lương tại ..." rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$data = New-Object 'object[,]' 37,2
$data[0,0] = 'Danh mục lương'
$data[0,1] = 4
$data[1,0] = 'Tổng công tại CẦN THƠ'
$data[1,1] = 0
$data[2,0] = 'Phụ cấp tại CẦN THƠ'
$data[2,1] = 0
$data[3,0] = 'Lương công tác tại CẦN THƠ'
$data[3,1] = 0
$data[4,0] = 'Lương cơ bản tại CẦN THƠ'
$data[4,1] = 0
$data[5,0] = 'Chiết khấu sale chính tại CẦN THƠ'
$data[5,1] = 0
$data[6,0] = 'Chiết khấu sale phụ tại CẦN THƠ'
$data[6,1] = 0
$data[7,0] = 'Đơn 1 bác sĩ tại CẦN THƠ'
$data[7,1] = 0
$data[8,0] = 'Đơn 2 bác sĩ tại CẦN THƠ'
$data[8,1] = 0
$data[9,0] = 'Công phụ phẫu 1 tại CẦN THƠ'
$data[9,1] = 0
$data[10,0] = 'Công phụ phẫu 2 tại CẦN THƠ'
$data[10,1] = 0
$data[11,0] = 'Ứng lương tại CẦN THƠ'
$data[11,1] = 0
$data[12,0] = 'Tổng công tại LONG XUYÊN'
$data[12,1] = 0
$data[13,0] = 'Phụ cấp tại LONG XUYÊN'
$data[13,1] = 0
$data[14,0] = 'Lương công tác tại LONG XUYÊN'
$data[14,1] = 0
$data[15,0] = 'Lương cơ bản tại LONG XUYÊN'
$data[15,1] = 0
$data[16,0] = 'Chiết khấu sale chính tại LONG XUYÊN'
$data[16,1] = 0
$data[17,0] = 'Chiết khấu sale phụ tại LONG XUYÊN'
$data[17,1] = 0
$data[18,0] = 'Đơn 1 bác sĩ tại LONG XUYÊN'
$data[18,1] = 0
$data[19,0] = 'Đơn 2 bác sĩ tại LONG XUYÊN'
$data[19,1] = 0
$data[20,0] = 'Công phụ phẫu 1 tại LONG XUYÊN'
$data[20,1] = 0
$data[21,0] = 'Công phụ phẫu 2 tại LONG XUYÊN'
$data[21,1] = 0
$data[22,0] = 'Ứng lương tại LONG XUYÊN'
$data[22,1] = 0
$data[23,0] = 'Tổng công tại SÓC TRĂNG'
$data[23,1] = 26
$data[24,0] = 'Phụ cấp tại SÓC TRĂNG'
$data[24,1] = 910000
$data[25,0] = 'Lương cơ bản tại SÓC TRĂNG'
$data[25,1] = 2785714.285714286
$data[26,0] = 'Chiết khấu sale chính tại SÓC TRĂNG'
$data[26,1] = 0
$data[27,0] = 'Chiết khấu sale phụ tại SÓC TRĂNG'
$data[27,1] = 0
$data[28,0] = 'Đơn 1 bác sĩ tại SÓC TRĂNG'
$data[28,1] = 0
$data[29,0] = 'Đơn 2 bác sĩ tại SÓC TRĂNG'
$data[29,1] = 0
$data[30,0] = 'Công phụ phẫu 1 tại SÓC TRĂNG'
$data[30,1] = 0
$data[31,0] = 'Công phụ phẫu 2 tại SÓC TRĂNG'
$data[31,1] = 0
$data[32,0] = 'Ứng lương tại SÓC TRĂNG'
$data[32,1] = 0
$data[33,0] = 'Tổng lương tại CẦN THƠ'
$data[33,1] = 0
$data[34,0] = 'Tổng lương tại LONG XUYÊN'
$data[34,1] = 0
$data[35,0] = 'Tổng lương tại SÓC TRĂNG'
$data[35,1] = 3695714.285714286
$data[36,0] = 'Tổng lương'
$data[36,1] = 3695714.285714286

$ws.Range("A1:B37").Value = $data

Write-Output "Updated Lương sheet rows 1-37"
